$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 5-7 (data now only spans rows 2-4 with updated TPM values)
$ws.Rows("5:7").Delete()

# Row 2: Sending=ECs, Target cluster ECs -> MuSCs, numeric columns recalculated with new TPM data
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.003549653112303053
$ws.Range("J2").Value = 0.003549653112303053
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.010147
$ws.Range("N2").Value = 3.030441
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.003455039455666667
$ws.Range("R2").Value = 0.031095355101
$ws.Range("S2").Value = 0.003549653112303053
$ws.Range("T2").Value = 0.003549653112303053

# Row 3: Sending cluster ECs -> FAPs, target cluster stays MuSCs, numeric columns recalculated
$ws.Range("A3").Value = "FAPs"
$ws.Range("I3").Value = 0.3907064193682856
$ws.Range("J3").Value = 0.3907064193682855
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.3802923980996668
$ws.Range("S3").Value = 0.3907064193682856
$ws.Range("T3").Value = 0.3907064193682855

# Row 4: Sending cluster FAPs -> MuSCs, target cluster ECs -> MuSCs, numeric columns recalculated
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5836756666666667
$ws.Range("H4").Value = 1.751027
$ws.Range("I4").Value = 0.6057439275194114
$ws.Range("J4").Value = 0.6057439275194113
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.010147
$ws.Range("N4").Value = 3.030441
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.5895982236563334
$ws.Range("R4").Value = 5.306384012907
$ws.Range("S4").Value = 0.6057439275194114
$ws.Range("T4").Value = 0.6057439275194113
